# Refresh the "cryptos" price/volume table (GitHub Actions scheduled update).
# Price (col D) and Volume(1h) (col E) cells are stored as literal text in the
# source sheet (e.g. "60.975.37", "  -2.10%  "), not numbers. Most new D
# values still contain a literal "." thousands separator so Excel keeps them
# as text automatically; a handful (e.g. "570.66") look like plain decimals,
# so for those we briefly force Text number-format before assigning the
# value and then restore the default "Normal" style so no stray formatting
# is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.964.21'
$ws.Range('E2').Value = '  -2.08%  '
$ws.Range('D3').Value = '2.421.45'
$ws.Range('E3').Value = '  -1.09%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.66'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.05'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.13%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  -1.14%  '
$ws.Range('D9').Value = '2.405.57'
$ws.Range('E9').Value = '  -1.54%  '
$ws.Range('E10').Value = '  -1.18%  '
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('E12').Value = '  -2.15%  '
$ws.Range('E13').Value = '  -1.82%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.19'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.03%  '
$ws.Range('E15').Value = '  -3.63%  '
$ws.Range('D16').Value = '2.831.62'
$ws.Range('E16').Value = '  -2.21%  '
$ws.Range('D17').Value = '60.909.25'
$ws.Range('E17').Value = '  -2.11%  '
$ws.Range('D18').Value = '2.392.16'
$ws.Range('E18').Value = '  -2.17%  '
$ws.Range('E19').Value = '  +6.76%  '
$ws.Range('E20').Value = '  -1.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.39%  '
$ws.Range('E22').Value = '  -1.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.07'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.37%  '
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('E25').Value = '  -2.79%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '64.81'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '585.41'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.17%  '
$ws.Range('E28').Value = '  -9.13%  '
$ws.Range('D29').Value = '2.537.25'
$ws.Range('E29').Value = '  -1.30%  '
$ws.Range('D30').Value = '0.0₃0935'
$ws.Range('E30').Value = '  -3.96%  '
$ws.Range('E32').Value = '  -4.33%  '
$ws.Range('E33').Value = '  -1.98%  '
$ws.Range('E34').Value = '  -2.67%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('E36').Value = '  -1.19%  '
$ws.Range('E37').Value = '  -5.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '151.85'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.98%  '
$ws.Range('E39').Value = '  -2.44%  '
$ws.Range('E40').Value = '  -0.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.15'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.13%  '
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('E43').Value = '  -1.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.14'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.07%  '
$ws.Range('E45').Value = '  -5.83%  '
$ws.Range('D46').Value = '0.0₆0291'
$ws.Range('E46').Value = '  +12.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '142.16'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.24%  '
$ws.Range('E48').Value = '  -3.42%  '
$ws.Range('E49').Value = '  -2.40%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.56'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.56%  '
$ws.Range('E51').Value = '  -3.38%  '

Write-Output "Applied 70 cell updates"
